$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.167.97"
$ws.Range("E2").Value = "  +0.97%  "
$ws.Range("D3").Value = "1.832.80"
$ws.Range("E3").Value = "  +0.87%  "
$ws.Range("E4").Value = "  +0.78%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.25"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.03%  "
$ws.Range("E6").Value = "  +0.79%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4706"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.00%  "
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07408"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8826"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.45"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.06%  "
$ws.Range("D12").Value = "1.821.42"
$ws.Range("E12").Value = "  -1.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07364"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +4.22%  "
$ws.Range("E14").Value = "  +2.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.87"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.556"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.011"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008789"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("E19").Value = "  +0.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.81"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.54%  "
$ws.Range("D21").Value = "27.193.58"
$ws.Range("E21").Value = "  +0.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.307"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.89%  "
$ws.Range("E23").Value = "  +0.99%  "
$ws.Range("D24").Value = "2.057.12"
$ws.Range("E24").Value = "  -1.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.901"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.50"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.58"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.159"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.280"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.74%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.53"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08936"
$ws.Range("D31").ClearFormats()
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7611"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.172"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.543"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.938"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.010"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.103"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05340"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01962"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.993"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.08%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.351"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.54%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.412"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.74%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5349"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1665"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.550"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4951"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("E47").Value = "  +1.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.009"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.672"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "103.86"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06319"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.52%  "
